$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("convnet")

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "imagenet"
$ws.Range("C13").Value = "u-net"
$ws.Range("D13").Value = "lab"
$ws.Range("E13").Value = 128
$ws.Range("F13").Value = 0.002
$ws.Range("G13").Value = "0.5 - #10"

$ws.Range("H13:O13").Select()
